$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet1 (nthRoot3): add Quo/Mod helper cells on rows 9 and 10
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("nthRoot3")
$ws1.Range("K9").Value = "Quo ="
$ws1.Range("L9").Value = 2
$ws1.Range("N9").Value = "Mod ="
$ws1.Range("O9").Value = 1
$ws1.Range("K10").Value = "Quo ="
$ws1.Range("L10").Value = 2
$ws1.Range("N10").Value = "Mod ="
$ws1.Range("O10").Value = 2

# ---------------------------------------------------------------------------
# 2) Sheet2 (nthRoot4): add Quo/mod helper cells on rows 5 and 6, plus new
#    "bundle 2" summary row 15
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("nthRoot4")
$ws2.Range("I5").Value = "Quo ="
$ws2.Range("J5").Value = 1
$ws2.Range("K5").Value = "mod ="
$ws2.Range("L5").Value = 3
$ws2.Range("I6").Value = "Quo ="
$ws2.Range("J6").Value = 2
$ws2.Range("K6").Value = "mod ="
$ws2.Range("L6").Value = 0
$ws2.Range("D15").Value = 7123
$ws2.Range("E15").Value = " = bundle 2"

# copy the visual style (big bold font) of the existing bundle-summary row
# onto the new row 15 without touching the values we just wrote
$ws2.Range("D9:F9").Copy()
$ws2.Range("D15:F15").PasteSpecial(-4122)
$ws2.Rows.Item(15).RowHeight = $ws2.Rows.Item(9).RowHeight

# ---------------------------------------------------------------------------
# 3) New worksheet nthRoot9-2, placed after nthRoot9
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("nthRoot9")
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws3)
$ws4.Name = "nthRoot9-2"

$ws4.Range("C3").Value = 9
$ws4.Range("D3").Value = "nth root"
$ws4.Range("C4").Value = 983271230
$ws4.Range("D4").Value = " bigInt"
$ws4.Range("C5").Value = 8
$ws4.Range("D5").Value = " magnitude"
$ws4.Range("C6").Value = 9
$ws4.Range("D6").Value = " number of digits"
$ws4.Range("C8").Value = "983271230/10^0"
$ws4.Range("D8").Value = "bundle 1"
$ws4.Range("C9").Formula = "=C4/POWER(10,0)"
$ws4.Range("D9").Value = "division"
$ws4.Range("C11").Value = 983271230
$ws4.Range("D11").Value = " = bundle 1"
$ws4.Range("C13").Value = 0
$ws4.Range("D13").Value = "Mod of Bundle 1"
$ws4.Range("C14").Value = ""
$ws4.Range("D14").Value = ""

$ws4.Columns.Item(3).ColumnWidth = 35.28515625

# copy styles (fonts/number formats) from the matching nthRoot3 cells so the
# new sheet reuses the existing style table instead of minting new entries
$ws1.Range("D7:E7").Copy()
$ws4.Range("C3:D3").PasteSpecial(-4122)
$ws1.Range("D8:E8").Copy()
$ws4.Range("C4:D4").PasteSpecial(-4122)
$ws1.Range("D9:E9").Copy()
$ws4.Range("C5:D5").PasteSpecial(-4122)
$ws1.Range("D10:E10").Copy()
$ws4.Range("C6:D6").PasteSpecial(-4122)
$ws1.Range("D12:E12").Copy()
$ws4.Range("C8:D8").PasteSpecial(-4122)
$ws1.Range("D13:E13").Copy()
$ws4.Range("C9:D9").PasteSpecial(-4122)
$ws1.Range("D17").Copy()
$ws4.Range("C11").PasteSpecial(-4122)
$ws1.Range("E15:F15").Copy()
$ws4.Range("D11:E11").PasteSpecial(-4122)
$ws1.Range("D21:E21").Copy()
$ws4.Range("C13:D13").PasteSpecial(-4122)
$ws1.Range("D9:E9").Copy()
$ws4.Range("C14:D14").PasteSpecial(-4122)

for ($r = 3; $r -le 14; $r++) {
    if ($r -eq 3 -or $r -eq 4 -or $r -eq 5 -or $r -eq 6 -or $r -eq 8 -or $r -eq 9 -or $r -eq 11 -or $r -eq 13 -or $r -eq 14) {
        $ws4.Rows.Item($r).RowHeight = $ws1.Rows.Item(7).RowHeight
    }
}

$ws4.Range("I10").Select()

# ---------------------------------------------------------------------------
# 4) View-state touch-ups (selection / scroll position / active sheet)
# ---------------------------------------------------------------------------
$ws1.Range("D8").Select()
$excel.ActiveWindow.ScrollRow = 20

$ws2.Range("M5").Select()
$excel.ActiveWindow.ScrollRow = 5

$ws3.Range("J10").Select()

$ws1.Activate()
$ws1.Range("D8").Select()
